# Update num_rows (D) and unique_rows (E) values for rows 2 and 3
# on the "dataInfo_test" sheet from 202176 to 213840.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataInfo_test")

$ws.Range("D2").Value = 213840
$ws.Range("E2").Value = 213840
$ws.Range("D3").Value = 213840
$ws.Range("E3").Value = 213840
